# Update odds/score values for the week of 2025-02-18 (FlashScore export)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.15
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 4.5
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.8
$ws.Range("AD2").Value = 9
$ws.Range("AG2").Value = 21
$ws.Range("AO2").Value = 17
$ws.Range("AR2").Value = 41

# Row 3
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93

# Row 4
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("U4").Value = 5
$ws.Range("V4").Value = 1.16

# Row 6
$ws.Range("G6").Value = 2.9
$ws.Range("I6").Value = 2.82
$ws.Range("J6").Value = 3.65
$ws.Range("K6").Value = 1.78
$ws.Range("L6").Value = 3.65
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 4.4
$ws.Range("O6").Value = 1.72
$ws.Range("P6").Value = 2
$ws.Range("S6").Value = 3.1
$ws.Range("T6").Value = 1.32
$ws.Range("W6").Value = 5.8
$ws.Range("X6").Value = 1.1
$ws.Range("Y6").Value = 1.7
$ws.Range("Z6").Value = 2.05
$ws.Range("AA6").Value = 2.4
$ws.Range("AB6").Value = 1.5
$ws.Range("AC6").Value = 5.8
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 12
$ws.Range("AF6").Value = 37
$ws.Range("AG6").Value = 37
$ws.Range("AH6").Value = 65
$ws.Range("AI6").Value = 4.4
$ws.Range("AJ6").Value = 5.6
$ws.Range("AK6").Value = 22
$ws.Range("AL6").Value = 175
$ws.Range("AN6").Value = 5.6
$ws.Range("AP6").Value = 12
$ws.Range("AQ6").Value = 37
$ws.Range("AR6").Value = 37
$ws.Range("AS6").Value = 70

# Row 7
$ws.Range("O7").Value = 1.73
$ws.Range("P7").Value = 2.1

# Row 8
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 10
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 3.25
$ws.Range("S8").Value = 2.08
$ws.Range("T8").Value = 1.73
$ws.Range("W8").Value = 3.75
$ws.Range("X8").Value = 1.25

# Row 9
$ws.Range("G9").Value = 2.38
$ws.Range("I9").Value = 2.9
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98
$ws.Range("AA9").Value = 2.05
$ws.Range("AB9").Value = 1.7
$ws.Range("AN9").Value = 7

# Row 11
$ws.Range("G11").Value = 1.85
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 4.5
$ws.Range("J11").Value = 2.63
$ws.Range("L11").Value = 5
$ws.Range("Y11").Value = 1.53
$ws.Range("Z11").Value = 2.38
$ws.Range("AG11").Value = 19
$ws.Range("AO11").Value = 21

# Row 12
$ws.Range("L12").Value = 4
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("Q12").Value = 1.78
$ws.Range("R12").Value = 2.1
$ws.Range("S12").Value = 2.35
$ws.Range("T12").Value = 1.57
$ws.Range("AC12").Value = 6.5
$ws.Range("AI12").Value = 7.5
$ws.Range("AP12").Value = 12
$ws.Range("AR12").Value = 29

# Row 14
$ws.Range("S14").Value = 2.25
$ws.Range("T14").Value = 1.62
$ws.Range("W14").Value = 4
$ws.Range("X14").Value = 1.22
$ws.Range("AM14").Value = 700
